$wb = $excel.ActiveWorkbook

# --- CourtMinuteOrders: selection moved to F2, no longer the active tab ---
$wsOrders = $wb.Worksheets.Item("CourtMinuteOrders")
$wsOrders.Activate()
$wsOrders.Range("F2").Select()

# --- CourtHearings: add two new rows of test data, becomes the active tab ---
$wsHearings = $wb.Worksheets.Item("CourtHearings")
$wsHearings.Activate()

$wsHearings.Range("A6").Value = "testT4149"
$wsHearings.Range("B6").Value = 1
$wsHearings.Range("C6").Value = 2
$wsHearings.Range("D6").Value = "Click"
$wsHearings.Range("E6").Value = "Yes"
$wsHearings.Range("F6").Value = "Contested"
$wsHearings.Range("G6").Value = "n/a"
$wsHearings.Range("H6").Value = "future"
$wsHearings.Range("I6").Value = "n/a"
$wsHearings.Range("J6").Value = "n/a"
$wsHearings.Range("K6").Value = "n/a"
$wsHearings.Range("L6").Value = "Click"
$wsHearings.Range("M6").Value = "n/a"
$wsHearings.Range("N6").Value = "n/a"
$wsHearings.Range("O6").Value = "n/a"
$wsHearings.Range("P6").Value = "n/a"
$wsHearings.Range("Q6").Value = "n/a"
$wsHearings.Range("R6").Value = "Yes"
$wsHearings.Range("S6").Value = "Yes"
$wsHearings.Range("T6").Value = "Yes"
$wsHearings.Range("U6").Value = "n/a"
$wsHearings.Range("V6").Value = "n/a"
$wsHearings.Range("W6").Value = "future"
$wsHearings.Range("X6").Value = "Auto"
$wsHearings.Range("Y6").Value = "Auto"
$wsHearings.Range("Z6").Value = "n/a"
$wsHearings.Range("AA6").Value = "n/a"
$wsHearings.Range("AB6").Value = "n/a"
$wsHearings.Range("AC6").Value = "n/a"
$wsHearings.Range("AD6").Value = "n/a"
$wsHearings.Range("A6:AD6").Interior.Color = 65535

$wsHearings.Range("A7").Value = "testT4149"
$wsHearings.Range("B7").Value = 1
$wsHearings.Range("C7").Value = 3
$wsHearings.Range("D7").Value = "Click"
$wsHearings.Range("E7").Value = "Yes"
$wsHearings.Range("F7").Value = "Continued"
$wsHearings.Range("G7").Value = "n/a"
$wsHearings.Range("H7").Value = "future"
$wsHearings.Range("I7").Value = "n/a"
$wsHearings.Range("J7").Value = "n/a"
$wsHearings.Range("K7").Value = "n/a"
$wsHearings.Range("L7").Value = "Click"
$wsHearings.Range("M7").Value = "n/a"
$wsHearings.Range("N7").Value = "n/a"
$wsHearings.Range("O7").Value = "n/a"
$wsHearings.Range("P7").Value = "n/a"
$wsHearings.Range("Q7").Value = "n/a"
$wsHearings.Range("R7").Value = "Yes"
$wsHearings.Range("S7").Value = "n/a"
$wsHearings.Range("T7").Value = "n/a"
$wsHearings.Range("U7").Value = "Yes"
$wsHearings.Range("V7").Value = "Yes"
$wsHearings.Range("W7").Value = "future"
$wsHearings.Range("X7").Value = "n/a"
$wsHearings.Range("Y7").Value = "n/a"
$wsHearings.Range("Z7").Value = "Courts Own Motion"
$wsHearings.Range("AA7").Value = "Auto"
$wsHearings.Range("AB7").Value = "n/a"
$wsHearings.Range("AC7").Value = "n/a"
$wsHearings.Range("AD7").Value = "n/a"
$wsHearings.Range("A7:AD7").Interior.Color = 65535

# Row 8: a single formatted (highlighted) but otherwise empty cell, F8
$wsHearings.Range("F8").Interior.Color = 65535

$wsHearings.Range("F10").Select()
